# Case_3_61 parallel contingency table: add two new columns (P, Q) and
# flip the I/K/M/O columns for rows 2-25, per the "contingencies with
# rene fine" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend the 0..13 header series with 14, 15 ---
# Copy the header style (bold/centered/bordered, same as O1) onto the
# two new header cells before writing their values.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows (2-25): swap the I/K and M/O column values ---
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# --- Data rows (2-25): new P and Q columns, all filled with 2 ---
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2

Write-Output "applied"
